# Insert a new record row right above the current row 142, pushing the
# existing rows 142:198 down to 143:199 (this mirrors the diff, where a new
# observation was inserted into the middle of the daily price log and every
# subsequent row's "shift" comes from that insertion).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(142).Insert()

$ws.Cells.Item(142, 1).Value = 10
$ws.Cells.Item(142, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(142, 3).Value = "La Araucanía"
$ws.Cells.Item(142, 4).Value = 44636
$ws.Cells.Item(142, 5).Value = 9
$ws.Cells.Item(142, 6).Value = 100112043
$ws.Cells.Item(142, 7).Value = "Pepino dulce"
$ws.Cells.Item(142, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(142, 9).Value = "Primera"
$ws.Cells.Item(142, 10).Value = 150
$ws.Cells.Item(142, 11).Value = 17000
$ws.Cells.Item(142, 12).Value = 17000
$ws.Cells.Item(142, 13).Value = 17000
$ws.Cells.Item(142, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(142, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(142, 16).Value = 944
$ws.Cells.Item(142, 17).Value = 18
$ws.Cells.Item(142, 18).Value = "Hortaliza"

# Make sure the Fecha column keeps the date/time number format used by the
# rest of the column (Rows.Insert should already carry it from row 141, but
# set it explicitly to be safe).
$ws.Cells.Item(142, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
